$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57
$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

# Column D holds dates stored as plain text (e.g. "2024-08-29"); force
# text entry so Excel doesn't auto-convert it to a date serial number,
# then reset the style back to the sheet default so no extra
# number-format style is left behind on the cell.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-08-29"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
